# Applies the PNAD 2009 "furto" data correction:
#  - Deletes two label-only rows (row 5 "situação do domicílio" and the row
#    that was originally row 8 "grandes regiões e unidades da federação"),
#    which shifts the numeric data for every region/category row up so it
#    lines up with the row immediately above its old position.
#  - Removes the two now-unused shared strings by virtue of the row deletes.
#  - Renames the "unnamed: 1_level_1" column header to "total".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 currently holds only the "situação do domicílio" label with no data.
# Deleting it shifts every row below (6..40) up by one.
$ws.Rows(5).Delete()

# After the first delete, the old row 8 ("grandes regiões e unidades da
# federação", also label-only) is now at row 7. Delete it too.
$ws.Rows(7).Delete()

# Row 2's "unnamed: 1_level_1" header becomes "total".
$ws.Range("B2").Value = "total"
